$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 14 from 45185 to 45204
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45185) {
        $cell.Value2 = 45204
    }
}
